# Insert a new weekly record at the top of the data table (row 9), pushing
# all existing rows down by one. The new row re-uses the constant columns
# (Mercado, Region, Categoria, etc.) from the row that follows it and gets
# its own date / volume / price values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 9:78 down to 10:79, creating a blank row at 9.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly observation.
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C9").Value = "Los Lagos"
$ws.Range("D9").Value = 45149
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = 100112012
$ws.Range("G9").Value = "Espinaca"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 13500
$ws.Range("N9").Value = "$/cuna 10 kilos"
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 1350
$ws.Range("Q9").Value = 10
$ws.Range("R9").Value = "Hortaliza"
